# Auto-generated Excel COM-interop edit script
# Applies the exact cell-level changes described by the target diff.
$wb = $excel.ActiveWorkbook


# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1951
$ws.Range("F4").Value = 109
$ws.Range("F5").Value = 453
$ws.Range("F8").Value = 1383
$ws.Range("F9").Value = 13
$ws.Range("F10").Value = 1058
$ws.Range("F11").Value = 1058
$ws.Range("F12").Value = 130
$ws.Range("F13").Value = 3017
$ws.Range("F15").Value = 895
$ws.Range("F16").Value = 1170
$ws.Range("F20").Value = 1728
$ws.Range("F22").Value = 1291
$ws.Range("F23").Value = 229
$ws.Range("F26").Value = 1088
$ws.Range("F27").Value = 1582
$ws.Range("F28").Value = 1489
$ws.Range("F30").Value = 397
$ws.Range("F31").Value = 1314
$ws.Range("F32").Value = 459
$ws.Range("F33").Value = 169
$ws.Range("F36").Value = 1876
$ws.Range("F37").Value = 495
$ws.Range("F41").Value = 2318
$ws.Range("F42").Value = 163
$ws.Range("F44").Value = 2834
$ws.Range("F46").Value = 846
$ws.Range("F47").Value = 651

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F13").Value = 378
$ws.Range("F14").Value = 0
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "2024-07-17"
$ws.Range("C20").Value = "上海·石川绫子小提琴动漫音乐会"
$ws.Range("D20").Value = "丁香路425号 上海东方艺术中心"
$ws.Range("E20").Value = "2024.07.17 19:30-07.17 21:00"
$ws.Range("F20").Value = 211
$ws.Range("G20").Value = 180
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=83967"
$ws.Range("I20").Value = "//i1.hdslb.com/bfs/openplatform/202404/HhY3CS7t1712652128640.jpeg"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "2024-07-19"
$ws.Range("C21").Value = "上海·《你的名字》《天气之子》《铃芽之旅》——新海诚动漫三部曲钢琴演奏会"
$ws.Range("D21").Value = "丁香路425号(上海科技馆地铁站1号口步行460米) 上海东方艺术中心音乐厅"
$ws.Range("E21").Value = "2024.07.19 19:30-07.19 21:30"
$ws.Range("F21").Value = 295
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=83479"
$ws.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202403/GpyueuYA1711508106584.jpeg"
$ws.Range("C22").Value = "上海·《时光机2008·青春重启计划》演唱会"
$ws.Range("D22").Value = "万航渡后路19号3楼 瓦肆VAS SHANGHAI"
$ws.Range("E22").Value = "2024.07.19 20:00-07.19 21:30"
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 128
$ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=86969"
$ws.Range("I22").Value = "//i0.hdslb.com/bfs/openplatform/202406/ajhKui2x1717588730159.png"
$ws.Range("C23").Value = "上海·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024"
$ws.Range("D23").Value = "宜昌路179号 万代南梦宫上海文化中心"
$ws.Range("E23").Value = "2024.07.19 19:30-07.19 21:00"
$ws.Range("F23").Value = 294
$ws.Range("G23").Value = 280
$ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=87061"
$ws.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202406/H9TGguhL1717747212685.png"
$ws.Range("F24").Value = 294
$ws.Range("F26").Value = 82
$ws.Range("F27").Value = 72
$ws.Range("F28").Value = 75
$ws.Range("F34").Value = 163
$ws.Range("F35").Value = 227
$ws.Range("F37").Value = 42
$ws.Range("F41").Value = 165

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F9").Value = 700
$ws.Range("F10").Value = 967
$ws.Range("F11").Value = 561
$ws.Range("F13").Value = 1397
$ws.Range("F15").Value = 1334

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 1951
$ws.Range("F4").Value = 700
$ws.Range("F5").Value = 967
$ws.Range("F6").Value = 561
$ws.Range("F7").Value = 109
$ws.Range("F8").Value = 1397
$ws.Range("F9").Value = 453
$ws.Range("F12").Value = 1383
$ws.Range("F13").Value = 13
$ws.Range("F14").Value = 1058
$ws.Range("F15").Value = 1058
$ws.Range("F16").Value = 3017
$ws.Range("F19").Value = 895
$ws.Range("F20").Value = 1170
$ws.Range("F23").Value = 1728
$ws.Range("F26").Value = 378
$ws.Range("F29").Value = 1088
$ws.Range("F30").Value = 1582
$ws.Range("F31").Value = 1489
$ws.Range("F33").Value = 397
$ws.Range("C34").Value = "上海·囚鸢代号鸢Only【女生专场】"
$ws.Range("D34").Value = "吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙"
$ws.Range("E34").Value = "2024.07.14 11:00-07.14 21:00"
$ws.Range("F34").Value = 1314
$ws.Range("G34").Value = 88
$ws.Range("H34").Value = "https://show.bilibili.com/platform/detail.html?id=87331"
$ws.Range("I34").Value = "//i1.hdslb.com/bfs/openplatform/202406/e053lkLC1718164512192.jpeg"
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "2024-07-19"
$ws.Range("C35").Value = "上海·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024"
$ws.Range("D35").Value = "宜昌路179号 万代南梦宫上海文化中心"
$ws.Range("E35").Value = "2024.07.19 19:30-07.19 21:00"
$ws.Range("F35").Value = 294
$ws.Range("G35").Value = 280
$ws.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=87061"
$ws.Range("I35").Value = "//i0.hdslb.com/bfs/openplatform/202406/H9TGguhL1717747212685.png"
$ws.Range("F36").Value = 459
$ws.Range("F39").Value = 1876
$ws.Range("F42").Value = 163
$ws.Range("F43").Value = 227
$ws.Range("F44").Value = 2318
$ws.Range("F45").Value = 163
$ws.Range("F47").Value = 2834
$ws.Range("F48").Value = 846
$ws.Range("F49").Value = 651

Write-Output "Applied all changes."